# Fix demo for TGOV1
# 1) Correct the misspelled sheet name "TOGV1DB" -> "TGOV1DB".
# 2) Update the saved selection on that sheet to cell L37 (last active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOGV1DB")
$ws.Name = "TGOV1DB"

$ws.Activate()
$ws.Range("L37").Select()
